# Add two new columns, I ("I0") and J ("IF"), to Sheet1.
# Mirrors the source diff: new header cells I1/J1 reuse the same header
# style as the existing header cells (e.g. H1), and new data cells I2:J72
# get plain (default) formatting like the existing data columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells (row 1) ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Reuse the exact header style/formatting already applied to H1 (bold,
# bordered, centered) by copying its format onto the two new header cells.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Data values (rows 2-72) ---
$iValues = @(3,9,7,4,7,12,8,9,7,8,9,8,7,8,8,6,8,7,6,9,8,5,7,8,9,7,9,7,7,5,5,7,8,9,8,7,6,8,8,7,8,8,8,9,7,9,8,6,7,10,7,7,7,9,11,7,9,8,7,8,7,8,8,8,6,6,5,3,3,4,3)
$jValues = @(4,9,7,5,8,12,8,9,7,8,9,8,7,8,8,7,8,7,7,9,8,5,7,8,9,7,9,8,7,5,5,7,8,9,8,7,6,8,8,7,8,8,8,9,7,9,8,6,7,10,7,7,7,9,11,7,9,8,7,8,7,9,8,8,6,6,5,3,3,4,3)

for ($r = 0; $r -lt $iValues.Length; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value  = $iValues[$r]
    $ws.Cells.Item($row, 10).Value = $jValues[$r]
}
